$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '35.335.43'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value2 = '  -0.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '1.909.82'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value2 = '  +2.12%  '
$ws.Range('B5').Value2 = 'XRP'
$ws.Range('C5').Value2 = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '0.695'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  +11.22%  '
$ws.Range('B6').Value2 = 'BNB'
$ws.Range('C6').Value2 = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '246.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '  +2.85%  '
$ws.Range('E7').Value2 = '  -0.45%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '41.81'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value2 = '  -1.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.349'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value2 = '  +5.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '52.50'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  +11.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '0.0726'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value2 = '  +3.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.0993'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value2 = '  +0.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '2.182.78'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value2 = '  +1.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '12.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value2 = '  +7.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '0.702'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '  +2.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '1.904.09'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value2 = '  +2.18%  '
$ws.Range('E17').Value2 = '  +1.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '35.307.19'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value2 = '  -0.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '72.59'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '  +2.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '0.0₃0820'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value2 = '  +2.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '240.51'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value2 = '  -1.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '12.50'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value2 = '  +1.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '4.84'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value2 = '  +1.78%  '
$ws.Range('E24').Value2 = '  -0.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '2.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value2 = '  +27.30%  '
$ws.Range('E26').Value2 = '  +0.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '170.50'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '  -0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '8.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value2 = '  +5.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '18.61'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value2 = '  +4.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '0.131'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value2 = '  +4.74%  '
$ws.Range('B32').Value2 = 'ImmutableX'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '0.978'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  +9.54%  '
$ws.Range('B33').Value2 = 'Filecoin'
$ws.Range('C33').Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '4.18'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value2 = '  +3.29%  '
$ws.Range('B34').Value2 = 'Hedera'
$ws.Range('C34').Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '0.0568'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value2 = '  +0.71%  '
$ws.Range('E35').Value2 = '  -0.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '4.11'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value2 = '  +1.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '1.75'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value2 = '  +2.40%  '
$ws.Range('E38').Value2 = '  -0.46%  '
$ws.Range('E39').Value2 = '  +1.30%  '
$ws.Range('B40').Value2 = 'Kaspa'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '0.0665'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value2 = '  +12.65%  '
$ws.Range('B41').Value2 = 'ARBITRUM'
$ws.Range('C41').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '1.11'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  +0.37%  '
$ws.Range('B42').Value2 = 'InjectiveProtocol'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '16.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value2 = '  +6.43%  '
$ws.Range('B43').Value2 = 'VeChain'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '0.0208'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value2 = '  +1.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '89.97'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value2 = '  -1.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '1.339.33'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value2 = '  -1.20%  '
$ws.Range('B46').Value2 = 'RenderToken'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '2.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value2 = '  +1.39%  '
$ws.Range('B47').Value2 = 'MultiversX'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '47.06'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value2 = '  -12.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '2.80'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '  +2.90%  '
$ws.Range('E49').Value2 = '  -0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '6.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value2 = '  -2.38%  '
$ws.Range('B51').Value2 = 'Gas'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '11.80'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value2 = '  -9.95%  '
